$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 449.3583473333333
$ws.Range("H2").Value = 1348.075042
$ws.Range("I2").Value = 0.959704436884883
$ws.Range("J2").Value = 0.9597044368848828
$ws.Range("M2").Value = 41.15160733333332
$ws.Range("N2").Value = 123.454822
$ws.Range("O2").Value = 0.511199339644863
$ws.Range("P2").Value = 0.511199339644863
$ws.Range("Q2").Value = 18491.81826141694
$ws.Range("R2").Value = 166426.3643527525
$ws.Range("S2").Value = 0.4906002743897973
$ws.Range("T2").Value = 0.4906002743897972
$ws.Range("G3").Value = 449.3583473333333
$ws.Range("H3").Value = 1348.075042
$ws.Range("I3").Value = 0.959704436884883
$ws.Range("J3").Value = 0.9597044368848828
$ws.Range("O3").Value = 0.02876767613861272
$ws.Range("P3").Value = 0.02876767613861272
$ws.Range("Q3").Value = 1040.624659898996
$ws.Range("R3").Value = 9365.621939090968
$ws.Range("S3").Value = 0.02760846642909401
$ws.Range("T3").Value = 0.027608466429094
$ws.Range("G4").Value = 449.3583473333333
$ws.Range("H4").Value = 1348.075042
$ws.Range("I4").Value = 0.959704436884883
$ws.Range("J4").Value = 0.9597044368848828
$ws.Range("M4").Value = 15.45528
$ws.Range("N4").Value = 46.36584
$ws.Range("O4").Value = 0.1919907736781588
$ws.Range("P4").Value = 0.1919907736781588
$ws.Range("Q4").Value = 6944.95907837392
$ws.Range("R4").Value = 62504.63170536527
$ws.Range("S4").Value = 0.1842543973398904
$ws.Range("T4").Value = 0.1842543973398904
$ws.Range("G5").Value = 449.3583473333333
$ws.Range("H5").Value = 1348.075042
$ws.Range("I5").Value = 0.959704436884883
$ws.Range("J5").Value = 0.9597044368848828
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.003444
$ws.Range("N5").Value = 0.010332
$ws.Range("O5").Value = 0.0000427825458061956
$ws.Range("P5").Value = 0.00004278254580619561
$ws.Range("Q5").Value = 1.547590148216
$ws.Range("R5").Value = 13.928311333944
$ws.Range("S5").Value = 0.00004105859903143666
$ws.Range("T5").Value = 0.00004105859903143666
$ws.Range("G6").Value = 449.3583473333333
$ws.Range("H6").Value = 1348.075042
$ws.Range("I6").Value = 0.959704436884883
$ws.Range("J6").Value = 0.9597044368848828
$ws.Range("M6").Value = 21.57398566666667
$ws.Range("N6").Value = 64.721957
$ws.Range("O6").Value = 0.2679994279925593
$ws.Range("P6").Value = 0.2679994279925593
$ws.Range("Q6").Value = 9694.450544566356
$ws.Range("R6").Value = 87250.05490109719
$ws.Range("S6").Value = 0.2572002401270699
$ws.Range("T6").Value = 0.2572002401270698
$ws.Range("I7").Value = 0.01202662913387072
$ws.Range("J7").Value = 0.01202662913387072
$ws.Range("M7").Value = 41.15160733333332
$ws.Range("N7").Value = 123.454822
$ws.Range("O7").Value = 0.511199339644863
$ws.Range("P7").Value = 0.511199339644863
$ws.Range("Q7").Value = 231.7320121629028
$ws.Range("R7").Value = 2085.588109466126
$ws.Range("S7").Value = 0.006148004871388381
$ws.Range("T7").Value = 0.00614800487138838
$ws.Range("I8").Value = 0.01202662913387072
$ws.Range("J8").Value = 0.01202662913387072
$ws.Range("O8").Value = 0.02876767613861272
$ws.Range("P8").Value = 0.02876767613861272
$ws.Range("S8").Value = 0.0003459781719623972
$ws.Range("T8").Value = 0.0003459781719623972
$ws.Range("I9").Value = 0.01202662913387072
$ws.Range("J9").Value = 0.01202662913387072
$ws.Range("M9").Value = 15.45528
$ws.Range("N9").Value = 46.36584
$ws.Range("O9").Value = 0.1919907736781588
$ws.Range("P9").Value = 0.1919907736781588
$ws.Range("Q9").Value = 87.03142756808001
$ws.Range("R9").Value = 783.28284811272
$ws.Range("S9").Value = 0.002309001832152123
$ws.Range("T9").Value = 0.002309001832152123
$ws.Range("I10").Value = 0.01202662913387072
$ws.Range("J10").Value = 0.01202662913387072
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.003444
$ws.Range("N10").Value = 0.010332
$ws.Range("O10").Value = 0.0000427825458061956
$ws.Range("P10").Value = 0.00004278254580619561
$ws.Range("Q10").Value = 0.019393775884
$ws.Range("R10").Value = 0.174543982956
$ws.Range("S10").Value = 0.0000005145298118139504
$ws.Range("T10").Value = 0.0000005145298118139505
$ws.Range("I11").Value = 0.01202662913387072
$ws.Range("J11").Value = 0.01202662913387072
$ws.Range("M11").Value = 21.57398566666667
$ws.Range("N11").Value = 64.721957
$ws.Range("O11").Value = 0.2679994279925593
$ws.Range("P11").Value = 0.2679994279925593
$ws.Range("Q11").Value = 121.4869462671201
$ws.Range("R11").Value = 1093.382516404081
$ws.Range("S11").Value = 0.003223129728556001
$ws.Range("T11").Value = 0.003223129728556
$ws.Range("G12").Value = 7.708291333333332
$ws.Range("H12").Value = 23.124874
$ws.Range("I12").Value = 0.01646276615823874
$ws.Range("J12").Value = 0.01646276615823874
$ws.Range("M12").Value = 41.15160733333332
$ws.Range("N12").Value = 123.454822
$ws.Range("O12").Value = 0.511199339644863
$ws.Range("P12").Value = 0.511199339644863
$ws.Range("Q12").Value = 317.2085781602697
$ws.Range("R12").Value = 2854.877203442427
$ws.Range("S12").Value = 0.008415755188819443
$ws.Range("T12").Value = 0.008415755188819441
$ws.Range("G13").Value = 7.708291333333332
$ws.Range("H13").Value = 23.124874
$ws.Range("I13").Value = 0.01646276615823874
$ws.Range("J13").Value = 0.01646276615823874
$ws.Range("O13").Value = 0.02876767613861272
$ws.Range("P13").Value = 0.02876767613861272
$ws.Range("Q13").Value = 17.85087134745511
$ws.Range("R13").Value = 160.657842127096
$ws.Range("S13").Value = 0.0004735955251859256
$ws.Range("T13").Value = 0.0004735955251859256
$ws.Range("G14").Value = 7.708291333333332
$ws.Range("H14").Value = 23.124874
$ws.Range("I14").Value = 0.01646276615823874
$ws.Range("J14").Value = 0.01646276615823874
$ws.Range("M14").Value = 15.45528
$ws.Range("N14").Value = 46.36584
$ws.Range("O14").Value = 0.1919907736781588
$ws.Range("P14").Value = 0.1919907736781588
$ws.Range("Q14").Value = 119.13380087824
$ws.Range("R14").Value = 1072.20420790416
$ws.Range("S14").Value = 0.003160699211602866
$ws.Range("T14").Value = 0.003160699211602865
$ws.Range("G15").Value = 7.708291333333332
$ws.Range("H15").Value = 23.124874
$ws.Range("I15").Value = 0.01646276615823874
$ws.Range("J15").Value = 0.01646276615823874
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.3333333333333333
$ws.Range("M15").Value = 0.003444
$ws.Range("N15").Value = 0.010332
$ws.Range("O15").Value = 0.0000427825458061956
$ws.Range("P15").Value = 0.00004278254580619561
$ws.Range("Q15").Value = 0.02654735535199999
$ws.Range("R15").Value = 0.238926198168
$ws.Range("S15").Value = 0.0000007043190472615357
$ws.Range("T15").Value = 0.0000007043190472615357
$ws.Range("G16").Value = 7.708291333333332
$ws.Range("H16").Value = 23.124874
$ws.Range("I16").Value = 0.01646276615823874
$ws.Range("J16").Value = 0.01646276615823874
$ws.Range("M16").Value = 21.57398566666667
$ws.Range("N16").Value = 64.721957
$ws.Range("O16").Value = 0.2679994279925593
$ws.Range("P16").Value = 0.2679994279925593
$ws.Range("Q16").Value = 166.2985667398242
$ws.Range("R16").Value = 1496.687100658418
$ws.Range("S16").Value = 0.004412011913583245
$ws.Range("T16").Value = 0.004412011913583244
$ws.Range("G17").Value = 1.356257333333333
$ws.Range("H17").Value = 4.068772
$ws.Range("I17").Value = 0.002896588408965574
$ws.Range("J17").Value = 0.002896588408965573
$ws.Range("M17").Value = 41.15160733333332
$ws.Range("N17").Value = 123.454822
$ws.Range("O17").Value = 0.511199339644863
$ws.Range("P17").Value = 0.511199339644863
$ws.Range("Q17").Value = 55.8121692242871
$ws.Range("R17").Value = 502.3095230185839
$ws.Range("S17").Value = 0.001480734081886166
$ws.Range("T17").Value = 0.001480734081886165
$ws.Range("G18").Value = 1.356257333333333
$ws.Range("H18").Value = 4.068772
$ws.Range("I18").Value = 0.002896588408965574
$ws.Range("J18").Value = 0.002896588408965573
$ws.Range("O18").Value = 0.02876767613861272
$ws.Range("P18").Value = 0.02876767613861272
$ws.Range("Q18").Value = 3.140822540876445
$ws.Range("R18").Value = 28.267402867888
$ws.Range("S18").Value = 0.00008332811725598113
$ws.Range("T18").Value = 0.00008332811725598111
$ws.Range("G19").Value = 1.356257333333333
$ws.Range("H19").Value = 4.068772
$ws.Range("I19").Value = 0.002896588408965574
$ws.Range("J19").Value = 0.002896588408965573
$ws.Range("M19").Value = 15.45528
$ws.Range("N19").Value = 46.36584
$ws.Range("O19").Value = 0.1919907736781588
$ws.Range("P19").Value = 0.1919907736781588
$ws.Range("Q19").Value = 20.96133683872
$ws.Range("R19").Value = 188.65203154848
$ws.Range("S19").Value = 0.0005561182496644876
$ws.Range("T19").Value = 0.0005561182496644874
$ws.Range("G20").Value = 1.356257333333333
$ws.Range("H20").Value = 4.068772
$ws.Range("I20").Value = 0.002896588408965574
$ws.Range("J20").Value = 0.002896588408965573
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 0.3333333333333333
$ws.Range("M20").Value = 0.003444
$ws.Range("N20").Value = 0.010332
$ws.Range("O20").Value = 0.0000427825458061956
$ws.Range("P20").Value = 0.00004278254580619561
$ws.Range("Q20").Value = 0.004670950256
$ws.Range("R20").Value = 0.042038552304
$ws.Range("S20").Value = 0.0000001239234262882649
$ws.Range("T20").Value = 0.0000001239234262882649
$ws.Range("G21").Value = 1.356257333333333
$ws.Range("H21").Value = 4.068772
$ws.Range("I21").Value = 0.002896588408965574
$ws.Range("J21").Value = 0.002896588408965573
$ws.Range("M21").Value = 21.57398566666667
$ws.Range("N21").Value = 64.721957
$ws.Range("O21").Value = 0.2679994279925593
$ws.Range("P21").Value = 0.2679994279925593
$ws.Range("Q21").Value = 29.2598762696449
$ws.Range("R21").Value = 263.338886426804
$ws.Range("S21").Value = 0.0007762840367326512
$ws.Range("T21").Value = 0.000776284036732651
$ws.Range("G22").Value = 4.171694666666667
$ws.Range("H22").Value = 12.515084
$ws.Range("I22").Value = 0.008909579414042005
$ws.Range("J22").Value = 0.008909579414042003
$ws.Range("M22").Value = 41.15160733333332
$ws.Range("N22").Value = 123.454822
$ws.Range("O22").Value = 0.511199339644863
$ws.Range("P22").Value = 0.511199339644863
$ws.Range("Q22").Value = 171.6719408372275
$ws.Range("R22").Value = 1545.047467535048
$ws.Range("S22").Value = 0.004554571112971738
$ws.Range("T22").Value = 0.004554571112971737
$ws.Range("G23").Value = 4.171694666666667
$ws.Range("H23").Value = 12.515084
$ws.Range("I23").Value = 0.008909579414042005
$ws.Range("J23").Value = 0.008909579414042003
$ws.Range("O23").Value = 0.02876767613861272
$ws.Range("P23").Value = 0.02876767613861272
$ws.Range("Q23").Value = 9.660816071326224
$ws.Range("R23").Value = 86.94734464193601
$ws.Range("S23").Value = 0.0002563078951144112
$ws.Range("T23").Value = 0.0002563078951144112
$ws.Range("G24").Value = 4.171694666666667
$ws.Range("H24").Value = 12.515084
$ws.Range("I24").Value = 0.008909579414042005
$ws.Range("J24").Value = 0.008909579414042003
$ws.Range("M24").Value = 15.45528
$ws.Range("N24").Value = 46.36584
$ws.Range("O24").Value = 0.1919907736781588
$ws.Range("P24").Value = 0.1919907736781588
$ws.Range("Q24").Value = 64.47470914784
$ws.Range("R24").Value = 580.27238233056
$ws.Range("S24").Value = 0.001710557044848921
$ws.Range("T24").Value = 0.001710557044848921
$ws.Range("G25").Value = 4.171694666666667
$ws.Range("H25").Value = 12.515084
$ws.Range("I25").Value = 0.008909579414042005
$ws.Range("J25").Value = 0.008909579414042003
$ws.Range("K25").Value = 1
$ws.Range("L25").Value = 0.3333333333333333
$ws.Range("M25").Value = 0.003444
$ws.Range("N25").Value = 0.010332
$ws.Range("O25").Value = 0.0000427825458061956
$ws.Range("P25").Value = 0.00004278254580619561
$ws.Range("Q25").Value = 0.014367316432
$ws.Range("R25").Value = 0.129305847888
$ws.Range("S25").Value = 0.0000003811744893951894
$ws.Range("T25").Value = 0.0000003811744893951894
$ws.Range("G26").Value = 4.171694666666667
$ws.Range("H26").Value = 12.515084
$ws.Range("I26").Value = 0.008909579414042005
$ws.Range("J26").Value = 0.008909579414042003
$ws.Range("M26").Value = 21.57398566666667
$ws.Range("N26").Value = 64.721957
$ws.Range("O26").Value = 0.2679994279925593
$ws.Range("P26").Value = 0.2679994279925593
$ws.Range("Q26").Value = 90.00008094437646
$ws.Range("R26").Value = 810.0007284993882
$ws.Range("S26").Value = 0.002387762186617538
$ws.Range("T26").Value = 0.002387762186617538
